# updated TSR for Consumer Staples & Industrials
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update TSR (column E) for the Industrials & Consumer Staples rows ---
$ws.Range("E61").Value = 58
$ws.Range("E62").Value = 225
$ws.Range("E63").Value = 34.72
$ws.Range("E64").Value = 121
$ws.Range("E65").Value = 45
$ws.Range("E66").Value = 112
$ws.Range("E67").Value = 211.93
$ws.Range("E68").Value = 167
$ws.Range("E69").Value = 84
$ws.Range("E70").Value = 87
$ws.Range("E71").Value = 123
$ws.Range("E72").Value = 38
$ws.Range("E73").Value = 101
$ws.Range("E74").Value = 84.32
$ws.Range("E75").Value = 40
$ws.Range("E76").Value = 222
$ws.Range("E77").Value = -51.22
$ws.Range("E78").Value = 59
$ws.Range("E79").Value = 38
$ws.Range("E80").Value = 40

# --- AutoFilter was toggled on the data range and then switched back off,
#     leaving behind the hidden _FilterDatabase defined name ---
$rng = $ws.Range("A1:E80")
$rng.AutoFilter()
$ws.AutoFilterMode = $false
$name = $ws.Names.Add("_xlnm._FilterDatabase", $rng)
$name.Visible = $false

# --- Selection moved to H16 ---
$ws.Range("H16").Select()
